$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 171 (shifts old rows 171-230 down to 172-231)
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new record
$ws.Cells.Item(171, 1).Value = 5
$ws.Cells.Item(171, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(171, 3).Value = "Maule"
$ws.Cells.Item(171, 4).NumberFormat = $ws.Cells.Item(172, 4).NumberFormat
$ws.Cells.Item(171, 4).Value = 44809
$ws.Cells.Item(171, 5).Value = 7
$ws.Cells.Item(171, 6).Value = 100112017
$ws.Cells.Item(171, 7).Value = "Apio"
$ws.Cells.Item(171, 8).Value = "Americana (o)"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 700
$ws.Cells.Item(171, 11).Value = 10000
$ws.Cells.Item(171, 12).Value = 10000
$ws.Cells.Item(171, 13).Value = 10000
$ws.Cells.Item(171, 14).Value = "$/docena de matas"
$ws.Cells.Item(171, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(171, 16).Value = 1667
$ws.Cells.Item(171, 17).Value = 6
$ws.Cells.Item(171, 18).Value = "Hortaliza"
